# HAJ-1: reflect reviewed Excel template changes for "Applicant Relatives Data"
# - Header cells D3/E3 (DateOfBirth / DateOfBirthHijri) gain an explicit
#   "DD/MM/YYYY" / "YYYYMMDD" format hint line.
# - Header cells H3/I3 (Relative DateOfBirth / Relative DateOfBirthHijri)
#   drop the redundant "القريب" word and gain the same format hint line.
# - The table column headers (bound to the same cells) follow automatically.
# - Active cell selection moves from H12 to F3.
# - Header/data columns are narrowed slightly to fit the new, taller header text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = "DateOfBirth`n`tتاريخ الميلاد`nDD/MM/YYYY"
$ws.Range("E3").Value = "DateOfBirthHijri`n`tتاريخ الميلاد هجري`nYYYYMMDD"
$ws.Range("H3").Value = "Relative DateOfBirth`n`tتاريخ الميلاد`nDD/MM/YYYY"
$ws.Range("I3").Value = "Relative DateOfBirthHijri`n`tتاريخ الميلاد هجري`nYYYYMMDD"

# Narrow columns B:J (closest the host's column-width quantisation allows to
# the reviewed template's 21.1796875 "characters" width).
$ws.Range("B1:J1").ColumnWidth = 20.3

# Reviewer's last selection was F3 (previously H12).
$ws.Range("F3").Select()

Write-Output "applicant-relatives-data header text, column width and selection updated"
